# Clarification Reflections - add two new discussion slides at the end
# (slide 15: "How large are these systems?" / slide 16: "How large is the
# clarifier for a community?") and bump the two stray date fields that
# live in the handout master / slide layout 13.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. "Last edited" date text fix living on the handout master.
# ---------------------------------------------------------------------
$handoutMaster = $p.HandoutMaster
for ($i = 1; $i -le $handoutMaster.Shapes.Count; $i++) {
    $shp = $handoutMaster.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "28/2/2023") {
            $tr.Text = "9/3/2023"
        }
    }
}

# Same literal date field also appears on the "1_Custom Layout" slide
# layout that belongs to the second slide master.
$design2 = $p.Designs.Item(2)
$layout13 = $design2.SlideMaster.CustomLayouts.Item(7)
for ($i = 1; $i -le $layout13.Shapes.Count; $i++) {
    $shp = $layout13.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "28/2/2023") {
            $tr.Text = "7/3/2023"
        }
    }
}

# ---------------------------------------------------------------------
# 2. New slide 15 - "How large are these systems?"
# ---------------------------------------------------------------------
$s15 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s15.Shapes.Item(1).TextFrame.TextRange.Text = "How large are these systems?"

$body15 = $s15.Shapes.Item(2).TextFrame.TextRange
$body15.Text = "What do you mean by " + [char]8220 + "large" + [char]8221 + "?`r" + `
    ("What do you mean by " + [char]8220 + "systems" + [char]8221) + "`r" + `
    "What is the context?`r" + `
    "How much does the size of these systems vary?`r" + `
    "What else do you need to know in order to estimate how large a system will be?"

# ---------------------------------------------------------------------
# 3. New slide 16 - "How large is the clarifier for a community?"
# ---------------------------------------------------------------------
$s16 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s16.Shapes.Item(1).TextFrame.TextRange.Text = "How large is the clarifier for a community?"

$body16 = $s16.Shapes.Item(2).TextFrame.TextRange
$body16.Text = "Define what else you need to know and make up answers for those additional inputs to your calculation.`r" + `
    "Rewrite the question so that it is clear and can be answered`r" + `
    "Then estimate the plan view area of the clarifier for the question you have crafted`r" + `
    "Hints:`r" + `
    "You can google to find out how much water people use (it is about 3 mL/s)`r" + `
    "Type your question (do this first) and your answer into the jam board`r" + `
    "`r" + `
    ""

$body16.Paragraphs(5, 1).IndentLevel = 2
$body16.Paragraphs(7, 1).IndentLevel = 2
